$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of A61 and A62
$a61 = $ws.Range("A61").Value()
$a62 = $ws.Range("A62").Value()
$ws.Range("A61").Value = $a62
$ws.Range("A62").Value = $a61

# Shorten the long source citation strings that are no longer needed:
# A66 previously held the long "STATISTICAL YEARBOOK..." citation; it now
# just repeats the short "State Statistical Committee" label (same as A65).
$ws.Range("A66").Value = $ws.Range("A65").Value()

# A68 previously held the long CESD citation; it now just repeats the
# short "CESD" label (same as A67).
$ws.Range("A68").Value = $ws.Range("A67").Value()
